# Add I0 and IF columns (I and J) to the worksheet, matching header row + data rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (the last existing header cell) onto the new
# header cells I1/J1 so they pick up the same bold/border/centered style
# used by the other headers, then set their text.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows: I{n} = I0 value, J{n} = IF value
$rows = @(
    @{ Row = 2; I = 1; J = 6 },
    @{ Row = 3; I = 1; J = 6 },
    @{ Row = 4; I = 1; J = 5 },
    @{ Row = 5; I = 1; J = 5 },
    @{ Row = 6; I = 1; J = 6 },
    @{ Row = 7; I = 1; J = 6 },
    @{ Row = 8; I = 1; J = 5 },
    @{ Row = 9; I = 1; J = 4 },
    @{ Row = 10; I = 1; J = 2 },
    @{ Row = 11; I = 7; J = 8 },
    @{ Row = 12; I = 8; J = 8 },
    @{ Row = 13; I = 6; J = 6 },
    @{ Row = 14; I = 7; J = 7 },
    @{ Row = 15; I = 6; J = 6 },
    @{ Row = 16; I = 7; J = 7 },
    @{ Row = 17; I = 7; J = 7 },
    @{ Row = 18; I = 7; J = 7 },
    @{ Row = 19; I = 8; J = 8 },
    @{ Row = 20; I = 5; J = 5 },
    @{ Row = 21; I = 7; J = 7 },
    @{ Row = 22; I = 7; J = 7 },
    @{ Row = 23; I = 8; J = 8 },
    @{ Row = 24; I = 7; J = 7 },
    @{ Row = 25; I = 11; J = 11 },
    @{ Row = 26; I = 8; J = 8 },
    @{ Row = 27; I = 7; J = 7 },
    @{ Row = 28; I = 5; J = 6 },
    @{ Row = 29; I = 9; J = 9 },
    @{ Row = 30; I = 4; J = 5 },
    @{ Row = 31; I = 6; J = 6 },
    @{ Row = 32; I = 5; J = 5 },
    @{ Row = 33; I = 7; J = 7 },
    @{ Row = 34; I = 7; J = 7 },
    @{ Row = 35; I = 7; J = 7 },
    @{ Row = 36; I = 7; J = 8 },
    @{ Row = 37; I = 5; J = 6 },
    @{ Row = 38; I = 6; J = 7 },
    @{ Row = 39; I = 3; J = 4 },
    @{ Row = 40; I = 5; J = 5 },
    @{ Row = 41; I = 6; J = 6 },
    @{ Row = 42; I = 5; J = 5 },
    @{ Row = 43; I = 6; J = 6 },
    @{ Row = 44; I = 4; J = 5 },
    @{ Row = 45; I = 1; J = 1 },
    @{ Row = 46; I = 6; J = 7 },
    @{ Row = 47; I = 6; J = 6 },
    @{ Row = 48; I = 5; J = 5 },
    @{ Row = 49; I = 7; J = 7 },
    @{ Row = 50; I = 5; J = 5 },
    @{ Row = 51; I = 4; J = 4 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 9).Value = $r.I
    $ws.Cells.Item($r.Row, 10).Value = $r.J
}
